$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add "Total" label in A17, bold (matching header style), and a SUM formula in C17
$ws.Range("A17").Value = "Total"
$ws.Range("A17").Font.Bold = $true

$ws.Range("C17").Formula = "=SUM(C3:C16)"

# Update selection to C12 as per the diff
$ws.Range("C12").Select()

# Set page orientation to portrait (triggers pageSetup element)
$ws.PageSetup.Orientation = 1

$wb.Save()
